{"js": "// Reorders the \"COMPETENCES TECHNIQUES\" skill lines.\n//\n// Before:\n//   Langages : ...\n//   Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\n//   Autres : anglais (confirm\u00e9)\n//   Visualisation : tableau\n//   ML/AI : ...\n//   MLOps : hadoop, spark, ...\n//\n// After:\n//   Langages : ...\n//   Visualisation : tableau\n//   MLOps : hadoop, spark, ...\n//   Autres : anglais (confirm\u00e9)\n//   ML/AI : ...\n//   Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Original text values that anchor the block we need to reorder (order-independent lookup).\nconst ORIG = {\n  langages: \"Langages : scala, r, javascript, python, matlab, c, c++\",\n  basesDeDonnees: \"Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\",\n  autres: \"Autres : anglais (confirm\u00e9)\",\n  visualisation: \"Visualisation : tableau\",\n  mlai: \"ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn\",\n  mlops: \"MLOps : hadoop, spark, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit\",\n};\n\n// Find the index (in document order) of each anchor paragraph.\nconst indexOf = {};\nfor (const key in ORIG) {\n  indexOf[key] = -1;\n}\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  for (const key in ORIG) {\n    if (t === ORIG[key]) {\n      indexOf[key] = i;\n    }\n  }\n}\n\nfor (const key in ORIG) {\n  if (indexOf[key] === -1) {\n    throw new Error(\"Could not locate paragraph for '\" + key + \"'\");\n  }\n}\n\n// The six paragraphs occupy six consecutive document positions; sort their\n// current indices so we can write the new text values into those same slots\n// in the right (new) order, regardless of original ordering.\nconst slots = Object.keys(indexOf)\n  .map((key) => indexOf[key])\n  .sort((a, b) => a - b);\n\nconst newOrder = [\n  ORIG.langages,\n  ORIG.visualisation,\n  ORIG.mlops,\n  ORIG.autres,\n  ORIG.mlai,\n  ORIG.basesDeDonnees,\n];\n\nfor (let i = 0; i < slots.length; i++) {\n  items[slots[i]].insertText(newOrder[i], \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Reorders the \"COMPETENCES TECHNIQUES\" skill lines.\n#\n# Before:\n#   Langages : ...\n#   Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\n#   Autres : anglais (confirm\u00e9)\n#   Visualisation : tableau\n#   ML/AI : ...\n#   MLOps : hadoop, spark, ...\n#\n# After:\n#   Langages : ...\n#   Visualisation : tableau\n#   MLOps : hadoop, spark, ...\n#   Autres : anglais (confirm\u00e9)\n#   ML/AI : ...\n#   Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\n\n$d = $word.ActiveDocument\n\n$langages       = \"Langages : scala, r, javascript, python, matlab, c, c++\"\n$basesDeDonnees = \"Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\"\n$autres         = \"Autres : anglais (confirm\u00e9)\"\n$visualisation  = \"Visualisation : tableau\"\n$mlai           = \"ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn\"\n$mlops          = \"MLOps : hadoop, spark, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit\"\n\n# Locate the six consecutive paragraphs by their current text (order-independent)\n# so the script is resilient to where exactly this block sits in the document.\n$slots = @{}\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd(\"`r\", \"`n\", [char]7)\n    if ($t -eq $langages)       { $slots[\"langages\"] = $i }\n    elseif ($t -eq $basesDeDonnees) { $slots[\"basesDeDonnees\"] = $i }\n    elseif ($t -eq $autres)     { $slots[\"autres\"] = $i }\n    elseif ($t -eq $visualisation) { $slots[\"visualisation\"] = $i }\n    elseif ($t -eq $mlai)       { $slots[\"mlai\"] = $i }\n    elseif ($t -eq $mlops)      { $slots[\"mlops\"] = $i }\n}\n\nforeach ($key in @(\"langages\", \"basesDeDonnees\", \"autres\", \"visualisation\", \"mlai\", \"mlops\")) {\n    if (-not $slots.ContainsKey($key)) {\n        throw \"Could not locate paragraph for '$key'\"\n    }\n}\n\n# Sort the six paragraph indices ascending, then write the new text values\n# into those positions in the desired final order.\n$indices = @($slots.Values | Sort-Object)\n\n$newOrder = @($visualisation, $mlops, $autres, $mlai, $basesDeDonnees)\n# index 0 (Langages paragraph) is left untouched; the remaining five indices\n# get the rest of the new ordering.\nfor ($k = 1; $k -lt $indices.Count; $k++) {\n    $d.Paragraphs.Item($indices[$k]).Range.Text = $newOrder[$k - 1]\n}\n\nWrite-Output \"done\"\n"}
